# Updates cryptocurrency market data (Price / Volume(1h) columns, and a few
# reordered coin rows) to match the latest scrape, mirroring the GitHub Actions
# commit "Updated cryptos list" workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.975.70'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.21'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.61'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5086'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2564'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06341'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07769'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.653.60'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.272'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5411'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.00'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7673'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.986.35'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '199.10'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.402'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.882'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.032'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.865'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.50'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1197'
$ws.Range('E26').Value = '  +4.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.812'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.61'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.234'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04902'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.252'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.164'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.523'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9087'
$ws.Range('E35').Value = '  +1.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.583'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.139.53'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5446'
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01563'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.525'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('B42').Value = 'BabyDogeCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₈126'
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8082'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.92'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.415'
$ws.Range('E45').Value = '  -4.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.777.22'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.86'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('E51').Value = '  -0.37%  '
